$d = $word.ActiveDocument

$replacements = @(
    @{old = "2025-09-30 Tuesday"; new = "2025-10-01 Wednesday"},
    @{old = "24×74=1776"; new = "91×83=7553"},
    @{old = "36×44=1584"; new = "76×14=1064"},
    @{old = "21×63=1323"; new = "75×13=975"},
    @{old = "45×14=630";  new = "92×76=6992"},
    @{old = "49×39=1911"; new = "93×53=4929"},
    @{old = "30×57=1710"; new = "34×59=2006"},
    @{old = "64×64=4096"; new = "96×44=4224"},
    @{old = "59×94=5546"; new = "19×54=1026"},
    @{old = "31×45=1395"; new = "47×13=611"},
    @{old = "48×11=528";  new = "85×13=1105"},
    @{old = "68×44=2992"; new = "75×25=1875"},
    @{old = "35×21=735";  new = "26×62=1612"},
    @{old = "39×89=3471"; new = "71×92=6532"},
    @{old = "79×47=3713"; new = "14×86=1204"},
    @{old = "17×51=867";  new = "82×48=3936"},
    @{old = "12×60=720";  new = "90×14=1260"},
    @{old = "72×58=4176"; new = "95×90=8550"},
    @{old = "99×61=6039"; new = "39×43=1677"},
    @{old = "44×61=2684"; new = "56×14=784"},
    @{old = "64×41=2624"; new = "63×47=2961"},
    @{old = "40×32=1280"; new = "60×39=2340"},
    @{old = "52×89=4628"; new = "42×24=1008"},
    @{old = "46×12=552";  new = "78×83=6474"},
    @{old = "15×42=630";  new = "97×70=6790"},
    @{old = "37×65=2405"; new = "90×28=2520"}
)

foreach ($r in $replacements) {
    $rng = $d.Content
    $rng.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}

$d.Save()
